$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row/column on the sheet
$lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11
$lastCol = $ws.Cells.SpecialCells(11).Column

# Old layout of the 6 metric groups (each 3 columns wide), starting at column F (6)
# Order: LLMBasedParaphrasing, InstructionPreservingCrossover, OnePointCrossover,
#        BertMLM, LLMBackTranslation_HI, BackTranslation_JA
# New layout:
#        BertMLM, OnePointCrossover, InstructionPreservingCrossover,
#        LLMBasedParaphrasing, BackTranslation_JA, LLMBackTranslation_HI

$groupStartCol = 6   # column F
$numGroups = 6
$groupWidth = 3

# Mapping from new group position (0-based) to old group position (0-based)
# new[0]=BertMLM(old idx3), new[1]=OnePointCrossover(old idx2), new[2]=InstructionPreservingCrossover(old idx1)
# new[3]=LLMBasedParaphrasing(old idx0), new[4]=BackTranslation_JA(old idx5), new[5]=LLMBackTranslation_HI(old idx4)
$oldIndexForNewGroup = @(3, 2, 1, 0, 5, 4)

# Snapshot the full range F1:W<lastRow> (18 columns x lastRow rows) before overwriting anything
$srcRange = $ws.Range($ws.Cells.Item(1, $groupStartCol), $ws.Cells.Item($lastRow, $groupStartCol + ($numGroups * $groupWidth) - 1))
$values = $srcRange.Value2

$totalCols = $numGroups * $groupWidth
$rowCount = $lastRow

# Build the new array with columns permuted per the group mapping.
# NOTE: New-Object 'object[,]' produces a 0-based .NET array, while the
# Value2 array returned from the Excel range is 1-based. Keep reads on the
# 1-based $values array and writes on the 0-based $newValues array.
$newValues = New-Object 'object[,]' $rowCount, $totalCols

for ($r = 1; $r -le $rowCount; $r++) {
    for ($g = 0; $g -lt $numGroups; $g++) {
        $oldG = $oldIndexForNewGroup[$g]
        for ($off = 0; $off -lt $groupWidth; $off++) {
            $newCol = ($g * $groupWidth) + $off + 1
            $oldCol = ($oldG * $groupWidth) + $off + 1
            $newValues[$r - 1, $newCol - 1] = $values[$r, $oldCol]
        }
    }
}

$destRange = $ws.Range($ws.Cells.Item(1, $groupStartCol), $ws.Cells.Item($lastRow, $groupStartCol + $totalCols - 1))
$destRange.Value2 = $newValues
